$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 277-278; this shifts the existing rows 277-319
# down to 279-321 (carrying their formatting, e.g. the date style on
# column D), matching the growth of the sheet's dimension from
# A1:T319 to A1:T321.
$ws.Rows("277:278").Insert()

# Row 277 (new): Ciruela / Blue Giant / Primera
$ws.Cells.Item(277,1).Value = 10
$ws.Cells.Item(277,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(277,3).Value = "La Araucanía"
$ws.Cells.Item(277,4).Value = 44984
$ws.Cells.Item(277,5).Value = 9
$ws.Cells.Item(277,6).Value = "Fruta"
$ws.Cells.Item(277,7).Value = 100103
$ws.Cells.Item(277,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(277,9).Value = 100103002
$ws.Cells.Item(277,10).Value = "Ciruela"
$ws.Cells.Item(277,11).Value = "Blue Giant"
$ws.Cells.Item(277,12).Value = "Primera"
$ws.Cells.Item(277,13).Value = 95
$ws.Cells.Item(277,14).Value = 15000
$ws.Cells.Item(277,15).Value = 15000
$ws.Cells.Item(277,16).Value = 15000
$ws.Cells.Item(277,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(277,18).Value = "Región de O'Higgins"
$ws.Cells.Item(277,19).Value = 833
$ws.Cells.Item(277,20).Value = 18

# Row 278 (new): Ciruela / Blue Giant / Primera (bins 450 kilos)
$ws.Cells.Item(278,1).Value = 10
$ws.Cells.Item(278,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(278,3).Value = "La Araucanía"
$ws.Cells.Item(278,4).Value = 44984
$ws.Cells.Item(278,5).Value = 9
$ws.Cells.Item(278,6).Value = "Fruta"
$ws.Cells.Item(278,7).Value = 100103
$ws.Cells.Item(278,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(278,9).Value = 100103002
$ws.Cells.Item(278,10).Value = "Ciruela"
$ws.Cells.Item(278,11).Value = "Blue Giant"
$ws.Cells.Item(278,12).Value = "Primera"
$ws.Cells.Item(278,13).Value = 8
$ws.Cells.Item(278,14).Value = 320000
$ws.Cells.Item(278,15).Value = 320000
$ws.Cells.Item(278,16).Value = 320000
$ws.Cells.Item(278,17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(278,18).Value = "Región de O'Higgins"
$ws.Cells.Item(278,19).Value = 711
$ws.Cells.Item(278,20).Value = 450
